$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Ensure Week_Start_Date column is treated as text so dates are stored as plain strings
$ws.Range("B2:B17").NumberFormat = "@"

$ws.Range("B2").Value = "2025-02-02"
$ws.Range("B3").Value = "2025-02-09"
$ws.Range("B4").Value = "2025-02-16"
$ws.Range("B5").Value = "2025-02-23"
$ws.Range("B6").Value = "2025-03-02"
$ws.Range("B7").Value = "2025-03-09"
$ws.Range("B8").Value = "2025-03-16"
$ws.Range("B9").Value = "2025-03-23"
$ws.Range("B10").Value = "2025-03-30"
$ws.Range("B11").Value = "2025-04-06"
$ws.Range("B12").Value = "2025-04-13"
$ws.Range("B13").Value = "2025-04-20"
$ws.Range("B14").Value = "2025-04-27"
$ws.Range("B15").Value = "2025-05-04"
$ws.Range("B16").Value = "2025-05-11"
$ws.Range("B17").Value = "2025-05-18"

# Inventory Coverage (L) updates
$ws.Range("L2").Value = 1.14
$ws.Range("L3").Value = 0.11
$ws.Range("L4:L17").Value = 0

# Stockout Risk / Reorder Urgency changes (row 2 only)
$ws.Range("M2").Value = "Low"
$ws.Range("N2").Value = "Normal"

# Seasonality Index (P) updates
$ws.Range("P2").Value = 0.86
$ws.Range("P3").Value = 1.04
$ws.Range("P4").Value = 1.07
$ws.Range("P5").Value = 0.91
$ws.Range("P6").Value = 0.87
$ws.Range("P7").Value = 0.88
$ws.Range("P8").Value = 0.87
$ws.Range("P9").Value = 0.87
$ws.Range("P10").Value = 0.85
$ws.Range("P11").Value = 1.16
$ws.Range("P12").Value = 0.83
$ws.Range("P13").Value = 1.18
$ws.Range("P14").Value = 0.91
$ws.Range("P15").Value = 0.95
$ws.Range("P16").Value = 1.13
$ws.Range("P17").Value = 1.04

# Remove "Sales Volume Rank" column (Q) and shift "Lifecycle Stage" (R) into its place, updating values to "Decline"
$ws.Range("Q1").Value = "Lifecycle Stage"
$ws.Range("Q2:Q17").Value = "Decline"

# Delete the now-duplicate Lifecycle Stage column R entirely
$ws.Range("R1:R17").ClearContents()

# Update Summary sheet: Max/Min Forecast Week become N/A (restock suggestion update)
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("B13").Value = "N/A"
$ws2.Range("B15").Value = "N/A"
